# Update the "修改时间" (last-modified timestamp) column on every portfolio
# sheet from 202509211606 -> 202509211628.
#
# The timestamp is stored as text (it looks numeric but must stay text, as
# in the source workbook), so each value is assigned with a leading
# apostrophe to force Excel to keep it as a text entry instead of coercing
# it to a number.

$wb = $excel.ActiveWorkbook

# old value being replaced everywhere: 202509211606
$newStamp = "202509211628"

# Sheet 1: 大智投资组合 - "修改时间" lives in column E, data rows 2-9.
$ws1 = $wb.Worksheets.Item(1)
for ($r = 2; $r -le 9; $r++) {
    $ws1.Range("E$r").Value = "'" + $newStamp
}

# Sheet 2: 大成投资组合 - "修改时间" lives in column E, data rows 2-11.
$ws2 = $wb.Worksheets.Item(2)
for ($r = 2; $r -le 11; $r++) {
    $ws2.Range("E$r").Value = "'" + $newStamp
}

# Sheet 3: 我的投资组合 - "修改时间" lives in column G, data rows 2-13.
$ws3 = $wb.Worksheets.Item(3)
for ($r = 2; $r -le 13; $r++) {
    $ws3.Range("G$r").Value = "'" + $newStamp
}
